$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 2, 1.02),
    @(2, 3, 1.023602512305199),
    @(2, 4, 1.028163376903938),
    @(2, 5, 1.034535632064959),
    @(2, 6, 1.046426151656681),
    @(2, 9, 1.031620878894135),
    @(2, 10, 1.028781809391155),
    @(2, 11, 1.030980993169081),
    @(2, 12, 1.037334826090936),
    @(2, 13, 1.049191621138973),
    @(2, 14, 1.013652721814535),
    @(3, 2, 1.02),
    @(3, 3, 1.024446788975063),
    @(3, 4, 1.028766109777419),
    @(3, 5, 1.03531672620725),
    @(3, 6, 1.047365262878193),
    @(3, 9, 1.031760979094724),
    @(3, 10, 1.029265180090462),
    @(3, 11, 1.031392157729511),
    @(3, 12, 1.037925203127292),
    @(3, 13, 1.049942041242463),
    @(3, 14, 1.013812939528151),
    @(4, 2, 1.02),
    @(4, 3, 1.024993492541244),
    @(4, 4, 1.029156046098251),
    @(4, 5, 1.035822885364214),
    @(4, 6, 1.047973737776813),
    @(4, 9, 1.031849929379007),
    @(4, 10, 1.029577715175837),
    @(4, 11, 1.031657442561358),
    @(4, 12, 1.038307298049421),
    @(4, 13, 1.050427794530968),
    @(4, 14, 1.013916510019741),
    @(5, 2, 1.02),
    @(5, 3, 1.025223421024438),
    @(5, 4, 1.029319956283175),
    @(5, 5, 1.036035849900419),
    @(5, 6, 1.048229732093692),
    @(5, 9, 1.031886915293336),
    @(5, 10, 1.029709046725102),
    @(5, 11, 1.031768783586002),
    @(5, 12, 1.038467948911229),
    @(5, 13, 1.050632047030354),
    @(5, 14, 1.013960026427327),
    @(6, 2, 1.02),
    @(6, 3, 1.025262032515439),
    @(6, 4, 1.029347476356061),
    @(6, 5, 1.03607161780915),
    @(6, 6, 1.048272725842531),
    @(6, 9, 1.031893101394525),
    @(6, 10, 1.029731094422134),
    @(6, 11, 1.031787467363934),
    @(6, 12, 1.038494923923075),
    @(6, 13, 1.050666344329721),
    @(6, 14, 1.013967331567039),
    @(7, 2, 1.02),
    @(7, 3, 1.024996564488199),
    @(7, 4, 1.029158236351234),
    @(7, 5, 1.035825730320385),
    @(7, 6, 1.047977157635276),
    @(7, 9, 1.031850425193942),
    @(7, 10, 1.029579470264631),
    @(7, 11, 1.031658931033044),
    @(7, 12, 1.038309444604464),
    @(7, 13, 1.050430523600549),
    @(7, 14, 1.013917091585416),
    @(8, 2, 1.02),
    @(8, 3, 1.023887755823251),
    @(8, 4, 1.028367087177364),
    @(8, 5, 1.034799452450808),
    @(8, 6, 1.046743360508216),
    @(8, 9, 1.03166857870066),
    @(8, 10, 1.028945215355049),
    @(8, 11, 1.031120106007078),
    @(8, 12, 1.037534329084999),
    @(8, 13, 1.049445190633059),
    @(8, 14, 1.01370688881402),
    @(9, 2, 1.02),
    @(9, 3, 1.021937023005544),
    @(9, 4, 1.026972506624288),
    @(9, 5, 1.03299675447827),
    @(9, 6, 1.044575511788682),
    @(9, 9, 1.031335130013882),
    @(9, 10, 1.027825806807181),
    @(9, 11, 1.030164817231098),
    @(9, 12, 1.036169160454163),
    @(9, 13, 1.047710366589786),
    @(9, 14, 1.013335729343562),
    @(10, 2, 1.02),
    @(10, 3, 1.020638726099384),
    @(10, 4, 1.026042571887354),
    @(10, 5, 1.031798911958201),
    @(10, 6, 1.043134591903545),
    @(10, 9, 1.031104133253634),
    @(10, 10, 1.027078415234001),
    @(10, 11, 1.029524128543068),
    @(10, 12, 1.035259587983246),
    @(10, 13, 1.046554891606232),
    @(10, 14, 1.01308780712498),
    @(11, 2, 1.02),
    @(11, 3, 1.020077087221802),
    @(11, 4, 1.025639871018728),
    @(11, 5, 1.031281192321378),
    @(11, 6, 1.042511702318238),
    @(11, 9, 1.031002056729334),
    @(11, 10, 1.026754535827392),
    @(11, 11, 1.029245811052098),
    @(11, 12, 1.034865877888992),
    @(11, 13, 1.046054832591447),
    @(11, 14, 1.012980344895862),
    @(12, 2, 1.02),
    @(12, 3, 1.019868550938002),
    @(12, 4, 1.025490286791547),
    @(12, 5, 1.031089033227099),
    @(12, 6, 1.042280491108593),
    @(12, 9, 1.030963833283274),
    @(12, 10, 1.026634195625243),
    @(12, 11, 1.029142298345366),
    @(12, 12, 1.034719658917124),
    @(12, 13, 1.045869130420726),
    @(12, 14, 1.012940412492541),
    @(13, 2, 1.02),
    @(13, 3, 1.01991327897206),
    @(13, 4, 1.025522373232384),
    @(13, 5, 1.031130245415084),
    @(13, 6, 1.042330079533106),
    @(13, 9, 1.030972046252232),
    @(13, 10, 1.026660010661165),
    @(13, 11, 1.029164508179444),
    @(13, 12, 1.034751022343943),
    @(13, 13, 1.045908962240605),
    @(13, 14, 1.012948978852539),
    @(14, 2, 1.019999999999999),
    @(14, 3, 1.020059847878631),
    @(14, 4, 1.025627506394844),
    @(14, 5, 1.03126530542236),
    @(14, 6, 1.042492587100133),
    @(14, 9, 1.030998903439202),
    @(14, 10, 1.026744589216104),
    @(14, 11, 1.029237257369325),
    @(14, 12, 1.034853790916722),
    @(14, 13, 1.046039481526049),
    @(14, 14, 1.012977044397906),
    @(15, 2, 1.02),
    @(15, 3, 1.020150164648959),
    @(15, 4, 1.025692282035376),
    @(15, 5, 1.031348539618822),
    @(15, 6, 1.042592734326398),
    @(15, 9, 1.031015410292402),
    @(15, 10, 1.026796695991616),
    @(15, 11, 1.029282062934178),
    @(15, 12, 1.034917113045853),
    @(15, 13, 1.046119904379074),
    @(15, 14, 1.012994334379655),
    @(16, 2, 1.02),
    @(16, 3, 1.020676011558596),
    @(16, 4, 1.026069297273245),
    @(16, 5, 1.031833291551173),
    @(16, 6, 1.043175953056293),
    @(16, 9, 1.03111086454365),
    @(16, 10, 1.027099904794402),
    @(16, 11, 1.029542580809301),
    @(16, 12, 1.035285720305828),
    @(16, 13, 1.046588084682433),
    @(16, 14, 1.013094936751816),
    @(17, 2, 1.02),
    @(17, 3, 1.021006004945771),
    @(17, 4, 1.026305781385028),
    @(17, 5, 1.032137620260039),
    @(17, 6, 1.04354206994616),
    @(17, 9, 1.031170191299659),
    @(17, 10, 1.027290032516661),
    @(17, 11, 1.0297057580574),
    @(17, 12, 1.035516976537961),
    @(17, 13, 1.046881835010518),
    @(17, 14, 1.013158012736164),
    @(18, 2, 1.02),
    @(18, 3, 1.021198535763489),
    @(18, 4, 1.026443715265521),
    @(18, 5, 1.032315221988681),
    @(18, 6, 1.043755719792136),
    @(18, 9, 1.03120459741427),
    @(18, 10, 1.027400906221499),
    @(18, 11, 1.029800850139865),
    @(18, 12, 1.035651877880217),
    @(18, 13, 1.047053200484203),
    @(18, 14, 1.013194793222848),
    @(19, 2, 1.02),
    @(19, 3, 1.021264192447397),
    @(19, 4, 1.026490746540626),
    @(19, 5, 1.032375795131908),
    @(19, 6, 1.043828585787267),
    @(19, 9, 1.031216295370705),
    @(19, 10, 1.027438707084696),
    @(19, 11, 1.029833259381071),
    @(19, 12, 1.035697878000018),
    @(19, 13, 1.04711163602408),
    @(19, 14, 1.013207332596055),
    @(20, 2, 1.02),
    @(20, 3, 1.020970594466093),
    @(20, 4, 1.026280409206585),
    @(20, 5, 1.032104959147544),
    @(20, 6, 1.04350277870673),
    @(20, 9, 1.031163846592786),
    @(20, 10, 1.027269636145188),
    @(20, 11, 1.029688259615228),
    @(20, 12, 1.035492163527373),
    @(20, 13, 1.046850315692481),
    @(20, 14, 1.01315124638037),
    @(21, 2, 1.02),
    @(21, 3, 1.020016684730131),
    @(21, 4, 1.025596547372937),
    @(21, 5, 1.031225529599392),
    @(21, 6, 1.042444728304376),
    @(21, 9, 1.030991003153936),
    @(21, 10, 1.02671968395083),
    @(21, 11, 1.029215838212851),
    @(21, 12, 1.034823527502014),
    @(21, 13, 1.046001045677353),
    @(21, 14, 1.012968780236666),
    @(22, 2, 1.02),
    @(22, 3, 1.019417395499079),
    @(22, 4, 1.025166558681532),
    @(22, 5, 1.030673437273108),
    @(22, 6, 1.041780403378545),
    @(22, 9, 1.030880549886287),
    @(22, 10, 1.026373693864588),
    @(22, 11, 1.028918038108807),
    @(22, 12, 1.034403260257181),
    @(22, 13, 1.045467319751684),
    @(22, 14, 1.012853963317733),
    @(23, 2, 1.02),
    @(23, 3, 1.019735044869397),
    @(23, 4, 1.025394504905581),
    @(23, 5, 1.030966031681399),
    @(23, 6, 1.042132487437129),
    @(23, 9, 1.030939271678687),
    @(23, 10, 1.026557129582314),
    @(23, 11, 1.029075980143061),
    @(23, 12, 1.034626039035449),
    @(23, 13, 1.045750234330921),
    @(23, 14, 1.012914838634207),
    @(24, 2, 1.02),
    @(24, 3, 1.020986594779402),
    @(24, 4, 1.026291873812298),
    @(24, 5, 1.032119717014511),
    @(24, 6, 1.043520532417577),
    @(24, 9, 1.031166714105139),
    @(24, 10, 1.027278852463579),
    @(24, 11, 1.029696166675277),
    @(24, 12, 1.035503375417245),
    @(24, 13, 1.046864557835777),
    @(24, 14, 1.013154303838308),
    @(25, 2, 1.02),
    @(25, 3, 1.022440954303498),
    @(25, 4, 1.027333084550734),
    @(25, 5, 1.033462105157555),
    @(25, 6, 1.045135200134708),
    @(25, 9, 1.031422871220043),
    @(25, 10, 1.02811540320976),
    @(25, 11, 1.030412462870842),
    @(25, 12, 1.036521999952118),
    @(25, 13, 1.048158677675268),
    @(25, 14, 1.013431769865446)
)

foreach ($row in $data) {
    $ws.Cells.Item($row[0], $row[1]).Value = $row[2]
}

Write-Output ("Updated " + $data.Length + " cells")
